$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 407, pushing existing rows 407-428 down to 409-430
$ws.Rows("407:408").Insert()

# New row 407: Conconina(o), week of 2021-11-08 (serial 44516)
$ws.Cells.Item(407, 1).Value = 7
$ws.Cells.Item(407, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(407, 3).Value = "Ñuble"
$ws.Cells.Item(407, 4).Value = 44516
$ws.Cells.Item(407, 5).Value = 16
$ws.Cells.Item(407, 6).Value = 100112033
$ws.Cells.Item(407, 7).Value = "Lechuga"
$ws.Cells.Item(407, 8).Value = "Conconina(o)"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 180
$ws.Cells.Item(407, 11).Value = 5000
$ws.Cells.Item(407, 12).Value = 5500
$ws.Cells.Item(407, 13).Value = 5250
$ws.Cells.Item(407, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(407, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(407, 16).Value = 525
$ws.Cells.Item(407, 17).Value = 10
$ws.Cells.Item(407, 18).Value = "Hortaliza"

# New row 408: Escarola, week of 2021-11-08 (serial 44516)
$ws.Cells.Item(408, 1).Value = 7
$ws.Cells.Item(408, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(408, 3).Value = "Ñuble"
$ws.Cells.Item(408, 4).Value = 44516
$ws.Cells.Item(408, 5).Value = 16
$ws.Cells.Item(408, 6).Value = 100112033
$ws.Cells.Item(408, 7).Value = "Lechuga"
$ws.Cells.Item(408, 8).Value = "Escarola"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 180
$ws.Cells.Item(408, 11).Value = 6500
$ws.Cells.Item(408, 12).Value = 7000
$ws.Cells.Item(408, 13).Value = 6750
$ws.Cells.Item(408, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(408, 15).Value = "Región del Maule"
$ws.Cells.Item(408, 16).Value = 450
$ws.Cells.Item(408, 17).Value = 15
$ws.Cells.Item(408, 18).Value = "Hortaliza"
